# Idiot proof features added to Report method
# Consolidate the report down to a single (gru3) row and bump its bandwidth.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the data for the row that should survive (previously row 7: the
# gru3 / google entry) before we start deleting rows. Read the numeric-
# looking text columns (C and I) via .Text so we keep their original
# "string" flavour instead of Excel's auto-numeric-coercion on .Value().
$survivorA = $ws.Range("A7").Value()
$survivorB = $ws.Range("B7").Value()
$survivorC = $ws.Range("C7").Text
$survivorD = $ws.Range("D7").Value()
$survivorE = $ws.Range("E7").Value()
$survivorF = $ws.Range("F7").Value()
$survivorG = $ws.Range("G7").Value()
$survivorH = $ws.Range("H7").Value()
$survivorI = $ws.Range("I7").Text

# Remove rows 3 through 7 (the extra data rows), leaving only the header
# (row 1) and the original row 2.
$ws.Range("A3:I7").EntireRow.Delete()

# Write the surviving row's data into row 2, replacing whatever was there.
$ws.Range("A2").Value = $survivorA
$ws.Range("B2").Value = $survivorB
$ws.Range("D2").Value = $survivorD
$ws.Range("E2").Value = $survivorE
$ws.Range("F2").Value = $survivorF
$ws.Range("G2").Value = $survivorG
$ws.Range("H2").Value = $survivorH

# Columns C and I hold small numeric-looking codes ("1", "2", ...) that are
# stored as text in this report, not as numbers. Force text formatting so
# the write doesn't get silently reinterpreted as a number, then restore
# the original (General) number format so no stray style is introduced.
$ws.Range("C2").NumberFormat = "@"
$ws.Range("C2").Value = $survivorC
$ws.Range("C2").NumberFormat = ""

$ws.Range("I2").NumberFormat = "@"
$ws.Range("I2").Value = $survivorI
$ws.Range("I2").NumberFormat = ""
